# Weekly update: insert 3 new "Coliflor" price rows (week of 2021-11-11)
# at the top of the date-ordered block that starts at row 401, pushing the
# existing rows 401:491 down to 404:494.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 401 (formatting/number-format of row 401, e.g. the
# date style on column D, is carried down onto the new rows automatically).
$ws.Rows("401:403").Insert()

$newRows = @(
    @(401, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", "2021-11-11", 13, 100112008, "Coliflor", "Sin especificar", "Primera", 12700, 500, 600, 554, "`$/unidad", "Región Metropolitana", 554, 1, "Hortaliza"),
    @(402, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", "2021-11-11", 13, 100112008, "Coliflor", "Sin especificar", "Primera", 3600,  700, 700, 700, "`$/unidad", "Región de Coquimbo",   700, 1, "Hortaliza"),
    @(403, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", "2021-11-11", 13, 100112008, "Coliflor", "Sin especificar", "Segunda", 3400,  400, 400, 400, "`$/unidad", "Región Metropolitana", 400, 1, "Hortaliza")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $col = 1
    for ($i = 1; $i -lt $r.Length; $i++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$i]
        $col = $col + 1
    }
}
